# Update the "enable WinRM" extension script URL used for the Windows rows
# (Windows Server 2012 R2 and Windows Server 2016 Datacenter) from the old
# ansible ConfigureRemotingForAnsible.ps1 script to Evgeny's new
# QualiSystems enable-winrm.ps1 script.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newScriptUrl = "https://raw.githubusercontent.com/QualiSystems/app-starter-pack/dev/scripts/enable-winrm.ps1"

$ws.Range("C2").Value = $newScriptUrl
$ws.Range("C3").Value = $newScriptUrl

# Move the view/selection: scroll so column B is the left-most visible
# column, and select C8.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 2
$ws.Range("C8").Select()
